# Start tracking whether recordings were sampled from a discourse or
# non-discourse context: populate the new "Discourse or standalone" (col B)
# values on the samples_retained sheet, and touch a few related notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# --- Column B: discourse / non-discourse classification per dataset row ---
$ws.Range("B2").Value  = "non-discourse"  # aesdd
$ws.Range("B4").Value  = "non-discourse"  # BAUM1
$ws.Range("B5").Value  = "discourse"      # BAUM2 [tr]
$ws.Range("B6").Value  = "discourse"      # BAUM2 [en]
$ws.Range("B7").Value  = "non-discourse"  # cafe
$ws.Range("B8").Value  = "non-discourse"  # CREMA-D
$ws.Range("B9").Value  = "non-discourse"  # dzafic
$ws.Range("B10").Value = "non-discourse"  # ekorpus
$ws.Range("B11").Value = "non-discourse"  # EmoDB
$ws.Range("B12").Value = "discourse"      # EmoReact_V_1.0
$ws.Range("B13").Value = "non-discourse"  # Emotional_EMA
$ws.Range("B14").Value = "non-discourse"  # EmoV-DB_sorted [en]
$ws.Range("B15").Value = "non-discourse"  # enterface_db
$ws.Range("B16").Value = "non-discourse"  # esd [en]
$ws.Range("B17").Value = "non-discourse"  # esd [zh]
$ws.Range("B18").Value = "discourse"      # EYASE
$ws.Range("B19").Value = "non-discourse"  # jl-corpus
$ws.Range("B20").Value = "discourse"      # LEGOv2
$ws.Range("B21").Value = "discourse"      # MELD
$ws.Range("B22").Value = "non-discourse"  # oreau2
$ws.Range("B23").Value = "non-discourse"  # ravdess
$ws.Range("B24").Value = "non-discourse"  # savee
$ws.Range("B25").Value = "discourse"      # ShEMO
$ws.Range("B26").Value = "non-discourse"  # tess

# --- Related note updates ---
# oreau2: new note clarifying versioning
$ws.Range("K22").Value = "not to be confused with v1"

# savee: extend the existing note about MetaData-derived speakers
$ws.Range("K24").Value = "added more unique speakers from MetaData folder; 4 main male speakers otherwise; the MetaData samples are arguably sampled from a discourse context"

# tess: clarify that the corpus is made of single words
$ws.Range("K26").Value = "Toronto English; single words"

# --- View state: scroll so row 2 is at top, with K27 selected ---
$ws.Activate()
$ws.Range("K27").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
